# Apply cryptos price/volume update (Fri Aug 25 15:58:19 UTC 2023 GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.011.85'
$ws.Range("E2").Value = '  -0.58%  '
$ws.Range("D3").Value = '1.651.25'
$ws.Range("E3").Value = '  -0.07%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.29%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.58'
$ws.Range("E5").Value = '  -0.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5210'
$ws.Range("E6").Value = '  +0.23%  '
$ws.Range("E8").Value = '  -1.14%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06269'
$ws.Range("E9").Value = '  -0.66%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.42'
$ws.Range("E10").Value = '  -3.61%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07740'
$ws.Range("E11").Value = '  +0.11%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.465'
$ws.Range("E12").Value = '  +1.02%  '
$ws.Range("D13").Value = '1.615.10'
$ws.Range("E13").Value = '  -2.21%  '
$ws.Range("D14").Value = '1.880.44'
$ws.Range("E14").Value = '  +0.08%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5424'
$ws.Range("E15").Value = '  -0.54%  '
$ws.Range("D16").Value = '0.0₅8089'
$ws.Range("E16").Value = '  -1.31%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.07'
$ws.Range("E17").Value = '  +0.60%  '
$ws.Range("D18").Value = '26.026.29'
$ws.Range("E18").Value = '  -0.60%  '
$ws.Range("E19").Value = '  -0.29%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.567'
$ws.Range("E20").Value = '  -2.41%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '191.34'
$ws.Range("E21").Value = '  -0.11%  '
$ws.Range("E22").Value = '  -1.40%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.976'
$ws.Range("E23").Value = '  -3.11%  '
$ws.Range("E24").Value = '  -0.36%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '138.73'
$ws.Range("E25").Value = '  +0.18%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1233'
$ws.Range("E26").Value = '  -0.52%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.254'
$ws.Range("E27").Value = '  -0.34%  '
$ws.Range("E28").Value = '  +0.43%  '
$ws.Range("E29").Value = '  -0.55%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05960'
$ws.Range("E30").Value = '  -1.60%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.275'
$ws.Range("E31").Value = '  -0.58%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.230'
$ws.Range("E33").Value = '  -3.78%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.547'
$ws.Range("E34").Value = '  -6.46%  '
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.412'
$ws.Range("E35").Value = '  -0.03%  '
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9452'
$ws.Range("E36").Value = '  -3.87%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.751'
$ws.Range("E37").Value = '  -0.63%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5698'
$ws.Range("E38").Value = '  -4.02%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01597'
$ws.Range("E39").Value = '  +0.05%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.884'
$ws.Range("E40").Value = '  -1.26%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8427'
$ws.Range("E41").Value = '  -2.39%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.003'
$ws.Range("E42").Value = '  -0.10%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.73'
$ws.Range("E43").Value = '  +1.04%  '
$ws.Range("D44").Value = '1.002.86'
$ws.Range("E44").Value = '  -4.75%  '
$ws.Range("D45").Value = '1.794.61'
$ws.Range("E45").Value = '  +0.05%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '56.71'
$ws.Range("E46").Value = '  -1.11%  '
$ws.Range("E47").Value = '  -4.87%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.9995'
$ws.Range("E48").Value = '  -0.54%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4299'
$ws.Range("E49").Value = '  +1.56%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.917'
$ws.Range("E50").Value = '  -1.71%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.483'
$ws.Range("E51").Value = '  +1.47%  '
